$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------------
# "No I'm running Windows XP. Can we use your laptop?" becomes three runs:
#   "No" + "," + " I'm running Windows XP. Can we use your laptop?"
# (i.e. a comma is inserted right after "No".)
$rng = $d.Content
$found = $rng.Find.Execute("No I’m running Windows XP. Can we use your laptop?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $rng.Start

    # Insert the comma right after "No" (collapsed range just after char 2).
    $insertPos = $d.Range($start + 2, $start + 2)
    $insertPos.InsertAfter(",")

    # The engine coalesces adjacent runs that share identical formatting, so
    # after the insert the paragraph is back to being one single <w:r>. Force
    # "No" and "," to stay in their own runs (as real Word would leave them
    # after a literal keystroke-by-keystroke edit) by toggling a character
    # property on and back off again right on those two sub-ranges; this is a
    # net-zero formatting change but it keeps the run boundary from being
    # re-merged with its neighbour.
    $runNo = $d.Range($start, $start + 2)
    $runNo.Bold = $true
    $runNo.Bold = $false

    $runComma = $d.Range($start + 2, $start + 3)
    $runComma.Bold = $true
    $runComma.Bold = $false

    Write-Output "Split 'No I’m running...' into 3 runs (No / , /  I’m running...)."
} else {
    Write-Output "WARNING: target sentence not found - no edit made."
}

# --- Change 2 ----------------------------------------------------------------
# Footer PAGE field's cached/displayed result goes from "5" to "1".
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
if ($footer.Range.Fields.Count -ge 1) {
    $pageField = $footer.Range.Fields.Item(1)
    Write-Output "Footer PAGE field result before: $($pageField.Result.Text)"
    $pageField.Result.Text = "1"
    Write-Output "Footer PAGE field result after: $($pageField.Result.Text)"
} else {
    Write-Output "WARNING: no PAGE field found in footer - no edit made."
}
